# Scheduled-runner data refresh: updates market-price derived columns
# (currentAveragePrice[NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ]) across
# several crafting-job sheets with freshly pulled values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 943.29034
$ws.Range("I15").Value = 943.29034
$ws.Range("K15").Value = 2829.87102
$ws.Range("M15").Value = -2660.87102
$ws.Range("H75").Value = 24289.25
$ws.Range("J75").Value = 27473.428
$ws.Range("L75").Value = 27473.428
$ws.Range("N75").Value = -29345.428
$ws.Range("H78").Value = 24289.25
$ws.Range("J78").Value = 27473.428
$ws.Range("L78").Value = 82420.284
$ws.Range("N78").Value = -91780.284
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3742.6667
$ws.Range("I32").Value = 3482.4558
$ws.Range("J32").Value = 5103.769
$ws.Range("K32").Value = 3482.4558
$ws.Range("L32").Value = 5103.769
$ws.Range("M32").Value = -3195.4558
$ws.Range("N32").Value = -5677.769
$ws.Range("H74").Value = 2864.225
$ws.Range("I74").Value = 2828.5151
$ws.Range("K74").Value = 2828.5151
$ws.Range("M74").Value = -1954.5151
$ws.Range("H77").Value = 2864.225
$ws.Range("I77").Value = 2828.5151
$ws.Range("K77").Value = 14142.5755
$ws.Range("M77").Value = -9774.575500000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1874.8684
$ws.Range("I134").Value = 1193.7037
$ws.Range("J134").Value = 3546.818
$ws.Range("K134").Value = 3581.1111
$ws.Range("L134").Value = 10640.454
$ws.Range("M134").Value = -1046.1111
$ws.Range("N134").Value = -15710.454
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7248742
$ws.Range("I31").Value = 1272.7755
$ws.Range("J31").Value = 25005042
$ws.Range("K31").Value = 1272.7755
$ws.Range("L31").Value = 25005042
$ws.Range("M31").Value = -977.7755
$ws.Range("N31").Value = -25005632
$ws.Range("H34").Value = 7248742
$ws.Range("I34").Value = 1272.7755
$ws.Range("J34").Value = 25005042
$ws.Range("K34").Value = 1272.7755
$ws.Range("L34").Value = 25005042
$ws.Range("M34").Value = -1070.7755
$ws.Range("N34").Value = -25005446
$ws.Range("H58").Value = 1652.4688
$ws.Range("I58").Value = 1571.5508
$ws.Range("J58").Value = 1859.2593
$ws.Range("K58").Value = 1571.5508
$ws.Range("L58").Value = 1859.2593
$ws.Range("M58").Value = -1368.5508
$ws.Range("N58").Value = -2265.2593
$ws.Range("H132").Value = 2164.6667
$ws.Range("I132").Value = 1754.2565
$ws.Range("J132").Value = 3943.111
$ws.Range("K132").Value = 5262.7695
$ws.Range("L132").Value = 11829.333
$ws.Range("M132").Value = -2732.7695
$ws.Range("N132").Value = -16889.333
$ws.Range("H136").Value = 1652.4688
$ws.Range("I136").Value = 1571.5508
$ws.Range("J136").Value = 1859.2593
$ws.Range("K136").Value = 4714.6524
$ws.Range("L136").Value = 5577.7779
$ws.Range("M136").Value = -2164.6524
$ws.Range("N136").Value = -10677.7779
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4953.846
$ws.Range("I3").Value = 3801.818
$ws.Range("J3").Value = 11290
$ws.Range("K3").Value = 11405.454
$ws.Range("L3").Value = 33870
$ws.Range("M3").Value = -11293.454
$ws.Range("N3").Value = -34094
$ws.Range("H107").Value = 26770.684
$ws.Range("I107").Value = 460.64285
$ws.Range("J107").Value = 42118.207
$ws.Range("K107").Value = 1381.92855
$ws.Range("L107").Value = 126354.621
$ws.Range("M107").Value = 538.0714499999999
$ws.Range("N107").Value = -130194.621
$ws.Range("H131").Value = 877.2461499999999
$ws.Range("J131").Value = 936.875
$ws.Range("L131").Value = 2810.625
$ws.Range("N131").Value = -12890.625
$ws.Range("H132").Value = 2312.0667
$ws.Range("I132").Value = 1216.6
$ws.Range("J132").Value = 2859.8
$ws.Range("K132").Value = 10949.4
$ws.Range("L132").Value = 25738.2
$ws.Range("M132").Value = -8419.4
$ws.Range("N132").Value = -30798.2
$ws.Range("H133").Value = 3554.7058
$ws.Range("I133").Value = 4732.857
$ws.Range("J133").Value = 2730
$ws.Range("K133").Value = 14198.571
$ws.Range("L133").Value = 8190
$ws.Range("M133").Value = -9138.571
$ws.Range("N133").Value = -18310
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 69504.5
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 69504.5
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 69504.5
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -70562.5
$ws.Range("H41").Value = 10919.429
$ws.Range("J41").Value = 21427.666
$ws.Range("L41").Value = 21427.666
$ws.Range("N41").Value = -22137.666
$ws.Range("H74").Value = 39245.25
$ws.Range("J74").Value = 39245.25
$ws.Range("L74").Value = 39245.25
$ws.Range("N74").Value = -41117.25
$ws.Range("H77").Value = 39245.25
$ws.Range("J77").Value = 39245.25
$ws.Range("L77").Value = 117735.75
$ws.Range("N77").Value = -127095.75
$ws.Range("H132").Value = 2254.3208
$ws.Range("I132").Value = 1461.4412
$ws.Range("J132").Value = 3673.158
$ws.Range("K132").Value = 4384.3236
$ws.Range("L132").Value = 11019.474
$ws.Range("M132").Value = -1854.3236
$ws.Range("N132").Value = -16079.474
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 13336.333
$ws.Range("I13").Value = 10000
$ws.Range("J13").Value = 15004.5
$ws.Range("K13").Value = 10000
$ws.Range("L13").Value = 15004.5
$ws.Range("M13").Value = -9860
$ws.Range("N13").Value = -15284.5
$ws.Range("H25").Value = 10451.5
$ws.Range("I25").Value = 5907
$ws.Range("J25").Value = 14996
$ws.Range("K25").Value = 5907
$ws.Range("L25").Value = 14996
$ws.Range("M25").Value = -5677
$ws.Range("N25").Value = -15456
$ws.Range("H26").Value = 11328.667
$ws.Range("J26").Value = 29997
$ws.Range("L26").Value = 29997
$ws.Range("N26").Value = -30587
$ws.Range("H45").Value = 22520.5
$ws.Range("I45").Value = 5041
$ws.Range("J45").Value = 40000
$ws.Range("K45").Value = 5041
$ws.Range("L45").Value = 40000
$ws.Range("M45").Value = -4634
$ws.Range("N45").Value = -40814
$ws.Range("H132").Value = 4125.579
$ws.Range("I132").Value = 1354.4445
$ws.Range("J132").Value = 8876.096
$ws.Range("K132").Value = 4063.3335
$ws.Range("L132").Value = 26628.288
$ws.Range("M132").Value = -1533.3335
$ws.Range("N132").Value = -31688.288
$ws.Range("H136").Value = 2516.4
$ws.Range("I136").Value = 1470.9429
$ws.Range("J136").Value = 6175.5
$ws.Range("K136").Value = 4412.8287
$ws.Range("L136").Value = 18526.5
$ws.Range("M136").Value = -1862.8287
$ws.Range("N136").Value = -23626.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 707.375
$ws.Range("I107").Value = 615.3333
$ws.Range("J107").Value = 825.7143
$ws.Range("K107").Value = 1845.9999
$ws.Range("L107").Value = 2477.1429
$ws.Range("M107").Value = 74.00009999999997
$ws.Range("N107").Value = -6317.1429
$ws.Range("H132").Value = 5556952.5
$ws.Range("I132").Value = 652.45
$ws.Range("J132").Value = 16669552
$ws.Range("K132").Value = 1957.35
$ws.Range("L132").Value = 50008656
$ws.Range("M132").Value = 572.6499999999999
$ws.Range("N132").Value = -50013716

Write-Host "Applied all changes"